$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 1.5
$ws.Range("I3").Value = 6.5
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("Z3").Value = 10
$ws.Range("AC3").Value = 8
$ws.Range("AE3").Value = 23
$ws.Range("AJ3").Value = 21
$ws.Range("AL3").Value = 51
$ws.Range("AN3").Value = 3.25
$ws.Range("AQ3").Value = 26
$ws.Range("AS3").Value = 201
$ws.Range("AU3").Value = 10
$ws.Range("AW3").Value = 8
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 151
$ws.Range("G4").Value = 2.82
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 2.45
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 3
$ws.Range("U4").Value = 1.7
$ws.Range("V4").Value = 1.93
$ws.Range("W4").Value = 9.25
$ws.Range("X4").Value = 15
$ws.Range("Y4").Value = 10
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 450
$ws.Range("AH4").Value = 7.8
$ws.Range("AI4").Value = 12
$ws.Range("AK4").Value = 27
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 30
$ws.Range("AO4").Value = 14.5
$ws.Range("AP4").Value = 19.5
$ws.Range("AQ4").Value = 60
$ws.Range("AR4").Value = 80
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.62
$ws.Range("AW4").Value = 4.4
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 19.5
$ws.Range("AZ4").Value = 55
$ws.Range("BA4").Value = 80
$ws.Range("BB4").Value = 200
